# adloori to davuluri completed
# Fill in the "Total Points" (column E) grading scores for the
# Customer Class and Product Class rubric sections - these scores
# were left blank before and now mirror the "Points for grading"
# values in column D for each graded item.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Customer Class section (rows 3-6)
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# Product Class section (rows 10-14)
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Leave the selection on the Product Class total (E15), matching
# where the grader's cursor ended up after entering the last score.
$ws.Range("E15").Select()
